{"js": "const body = context.document.body;\n\n// 1) Remove the \"To be provided separately as a word doc for students to\n//    include with every submission\" paragraph that sits under the\n//    \"Assessment Cover Page\" heading.\nconst allParagraphs = body.paragraphs;\nallParagraphs.load(\"items/text\");\nawait context.sync();\n\nlet noteParagraph = null;\nfor (const p of allParagraphs.items) {\n  if (p.text.indexOf(\"To be provided separately\") !== -1) {\n    noteParagraph = p;\n    break;\n  }\n}\nif (noteParagraph) {\n  noteParagraph.delete();\n}\n\n// 2) Find the \"Assessment Due Date\" and \"Date of Submission\" rows in the\n//    cover-sheet table so the edits aren't tied to brittle fixed indices.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst coverTable = tables.items[0];\nconst rows = coverTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet dueDateValueCell = null;\nlet submissionValueCell = null;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  const labelCell = cells.items[0];\n  labelCell.body.load(\"text\");\n  await context.sync();\n\n  const label = labelCell.body.text.trim();\n  if (label.indexOf(\"Assessment Due Date\") === 0) {\n    dueDateValueCell = cells.items[1];\n  } else if (label.indexOf(\"Date of Submission\") === 0) {\n    submissionValueCell = cells.items[1];\n  }\n}\n\n// 3) Add the extension note as a brand-new paragraph right after\n//    \"26 May 2023\" in the \"Assessment Due Date\" cell. Using insertText with\n//    an embedded paragraph mark (\\r) \u2014 rather than Paragraph.insertParagraph \u2014\n//    keeps the new paragraph's formatting plain/default instead of it\n//    picking up unrelated heading formatting.\nif (dueDateValueCell) {\n  const dueDateParagraphs = dueDateValueCell.body.paragraphs;\n  dueDateParagraphs.load(\"items/text\");\n  await context.sync();\n\n  const lastDueDatePara = dueDateParagraphs.items[dueDateParagraphs.items.length - 1];\n  const endOfLastPara = lastDueDatePara.getRange(\"End\");\n  endOfLastPara.insertText(\n    \"\\rExtension Granted to 16 June 2023 due to Personal Mitigating Circumstances\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n\n// 4) Fill in the previously empty \"Date of Submission\" value cell.\nif (submissionValueCell) {\n  const submissionParagraphs = submissionValueCell.body.paragraphs;\n  submissionParagraphs.load(\"items/text\");\n  await context.sync();\n\n  const lastSubmissionPara = submissionParagraphs.items[submissionParagraphs.items.length - 1];\n  lastSubmissionPara.insertText(\"16 June 2023\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the \"To be provided separately as a word doc for students to\n#    include with every submission\" paragraph right under \"Assessment Cover Page\".\n$target = $d.Paragraphs.Item(4)\nif ($target.Range.Text -like \"To be provided separately*\") {\n    $target.Range.Delete()\n}\n\n# 2) Add the extension note as a new paragraph after \"26 May 2023\" in the\n#    \"Assessment Due Date\" cell (table 1, row 6, column 2).\n$dueCell = $d.Tables.Item(1).Cell(6, 2)\n$dueRange = $dueCell.Range\n$dueRange.Collapse(0) | Out-Null         # wdCollapseEnd\n$dueRange.MoveEnd(1, -1) | Out-Null      # step back before the end-of-cell marker\n$dueRange.InsertAfter(\"`r\" + \"Extension Granted to 16 June 2023 due to Personal Mitigating Circumstances\") | Out-Null\n\n# 3) Fill in the submission date in the previously empty \"Date of Submission\"\n#    cell (table 1, row 7, column 2).\n$subCell = $d.Tables.Item(1).Cell(7, 2)\n$subRange = $subCell.Range\n$subRange.Collapse(0) | Out-Null         # wdCollapseEnd\n$subRange.MoveEnd(1, -1) | Out-Null      # step back before the end-of-cell marker\n$subRange.Text = \"16 June 2023\"\n"}
